# Update 3-Mar-2021, end of day.
# Fill in the "Buku KAS HARIAN" daily transactions on Sheet1 for
# 1-Mar-2021 .. 3-Mar-2021 (rows 3-22), matching the paper cash book.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1 Mar 2021 (row 3 already has the date + "Wages Expense") ---
$ws.Range("D3").Formula  = "=45000+195000"

$ws.Range("B4").Value    = "A/R"
$ws.Range("C4").Formula  = "=5400000+810000+43500000+45874000"

$ws.Range("B5").Value    = "TRANSFER BCA"
$ws.Range("D5").Formula  = "=1800000+2567000+810000+30000000+720000"

$ws.Range("B6").Value    = "SALES - cash/retail"
$ws.Range("C6").Formula  = "=69749025-14455025-45874000"

$ws.Range("B7").Value    = "SELISIH - kurang"
$ws.Range("D7").Value    = 50000

$ws.Range("B8").Value    = "SETOR KE BANK"
$ws.Range("D8").Value    = 69000000

# --- 2 Mar 2021 ---
$ws.Range("A9").Value    = 44257
$ws.Range("B9").Value    = "Wages Expense"
$ws.Range("D9").Formula  = "=45000+180000"

$ws.Range("B10").Value   = "TRANSFER BCA"
$ws.Range("D10").Formula = "=900000+1519000+220000"

$ws.Range("B11").Value   = "BELI kresek"
$ws.Range("D11").Formula = "=52000"

$ws.Range("B12").Value   = "A/R"
$ws.Range("C12").Formula = "=9000000+9566000"

$ws.Range("B13").Value   = "SALES - cash/retail"
$ws.Range("C13").Formula = "=2231975+16952025-9566000"

$ws.Range("B14").Value   = "SELISIH - lebih"
$ws.Range("C14").Value   = 470000

$ws.Range("B15").Value   = "SETOR KE BANK"
$ws.Range("D15").Value   = 26000000

# --- 3 Mar 2021 ---
$ws.Range("A16").Value   = 44258
$ws.Range("B16").Value   = "Wages Expense"
$ws.Range("D16").Formula = "=45000+210000"

$ws.Range("B17").Value   = "A/R"
$ws.Range("C17").Formula = "=1744500+7700000+84925000+18837500"

$ws.Range("B18").Value   = "TRANSFER BCA"
$ws.Range("D18").Formula = "=1744500+7700000+7057500+84925000+1365000"

$ws.Range("B19").Value   = "SALES - cash/retail"
$ws.Range("C19").Formula = "=8225475+18785525-18837500"

$ws.Range("B20").Value   = "SELISIH - kurang"
$ws.Range("C20").Value   = " "
$ws.Range("D20").Value   = 100000

$ws.Range("B21").Value   = "SETOR KE BANK"
$ws.Range("D21").Value   = 18000000

# --- 4 Mar 2021 starts, only the date + label recorded so far ---
$ws.Range("A22").Value   = 44259
$ws.Range("B22").Value   = "Wages Expense"

# Put the cursor where the bookkeeper left off and scroll the frozen
# pane down so row 20 is visible at the top of the data area.
$ws.Range("C20").Select()
$app = $excel
$win = $app.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
